# Add a default header containing the questionnaire number ("Questionnaire 12")
# to the document's single section, so printed copies can be matched back to
# the originating questionnaire.

$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1 -- the "default" header for the section.
$header = $d.Sections(1).Headers(1)

# Insert the header text. Using InsertAfter (rather than assigning .Text)
# keeps Word from also materializing the first-page/even-page headers and
# their paired footers, matching a minimal "just add a default header" edit.
$header.Range.InsertAfter("Questionnaire 12")

# Style + center the paragraph (matches the Header style, centered).
$headerRange = $header.Range
$headerRange.Style = "Header"
$headerRange.ParagraphFormat.Alignment = 1

# Apply the run formatting (Arial, 12pt) only to the visible text, not the
# trailing paragraph mark, so the paragraph mark doesn't pick up its own
# run-properties block.
$textRange = $headerRange.Duplicate
$textRange.End = $headerRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
